# Generate Report for Handoff
# b.md has now been handed off (target xliffs generated), so update the
# status/tracking rows for "b.md" (row 3) on all three report sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: b.md row (row 3) ------------------------------------
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-26 14:48:31"

# --- zh-cn sheet: b.md row (row 3) ----------------------------------------
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces this to stay a text value ("False") instead of
# being auto-coerced to the boolean FALSE; reset the style afterwards so no
# quote-prefix formatting sticks to the cell.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-26 14:48:26"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b58606a87955b29669492dd45638ad7614be31d/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3226d6627a87edcac68b62d4e92b1c9f19c82e3/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.1667

# --- de-de sheet: b.md row (row 3) ----------------------------------------
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-26 14:48:31"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b58606a87955b29669492dd45638ad7614be31d/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3226d6627a87edcac68b62d4e92b1c9f19c82e3/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.1667
